# Update the "想去人数" (F column) figures on the "展览" and "全部类型"
# sheets to reflect the newly scraped counts.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (rows are offset by one compared to "全部类型" because that
# sheet has an extra data row that "展览" does not).
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F2").Value = 1325
$wsExhibit.Range("F3").Value = 1214
$wsExhibit.Range("F4").Value = 14589
$wsExhibit.Range("F5").Value = 17611
$wsExhibit.Range("F23").Value = 213
$wsExhibit.Range("F24").Value = 7244
$wsExhibit.Range("F25").Value = 980
$wsExhibit.Range("F30").Value = 5861
$wsExhibit.Range("F33").Value = 139
$wsExhibit.Range("F36").Value = 5083
$wsExhibit.Range("F37").Value = 21

# Sheet "全部类型"
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value = 1325
$wsAll.Range("F3").Value = 1214
$wsAll.Range("F4").Value = 14589
$wsAll.Range("F5").Value = 17611
$wsAll.Range("F24").Value = 213
$wsAll.Range("F25").Value = 7244
$wsAll.Range("F26").Value = 980
$wsAll.Range("F32").Value = 5861
$wsAll.Range("F35").Value = 139
$wsAll.Range("F38").Value = 5083
$wsAll.Range("F39").Value = 21
